# NSMB - Several more rooms of 8-2 done
#
# Adds timing data for several newly-completed rooms of World 8-2:
#  - New split timers in I93/J93 for the "Get flag" row
#  - Five new rows (102-106) of room timings below the existing table
#  - Updates the active selection to J94

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New split values recorded alongside the existing "Get flag" row (row 93)
$ws.Range("I93").Value = 31543
$ws.Range("J93").Value = 36872

# New rooms appended to the 8-2 timing table
$ws.Range("A102").Value = "Rail 93650944 (after turn)"
$ws.Range("B102").Value = 30698
$ws.Range("C102").Value = 35992

$ws.Range("A103").Value = "Black screens"
$ws.Range("B103").Value = 30992
$ws.Range("C103").Value = 36290

$ws.Range("A104").Value = "Black screens"
$ws.Range("B104").Value = 31261
$ws.Range("C104").Value = 36576

$ws.Range("A105").Value = "Black screen"
$ws.Range("B105").Value = 31543
$ws.Range("C105").Value = 36872

$ws.Range("A106").Value = "Black screen (water scene)"
$ws.Range("B106").Value = 31909
$ws.Range("C106").Value = 37235

# Extend the existing "elapsed" formula down through the new rows
# (set per-cell rather than as one range assignment so each new formula
# cell carries its own formula text instead of an incomplete shared-group
# stub without a master)
$ws.Range("D102").Formula = "=IF(B102 >  0,C102-B102, 0)"
$ws.Range("D103").Formula = "=IF(B103 >  0,C103-B103, 0)"
$ws.Range("D104").Formula = "=IF(B104 >  0,C104-B104, 0)"
$ws.Range("D105").Formula = "=IF(B105 >  0,C105-B105, 0)"
$ws.Range("D106").Formula = "=IF(B106 >  0,C106-B106, 0)"

# Match the author's final cursor position
$ws.Range("J94").Select()
